# Fruta / hortaliza, semanal
# Insert two new weekly observation rows (kiwi, Feria Lagunitas de Puerto
# Montt) above the current row 120, pushing the existing data down by two
# rows (old row 120 -> 122, ..., old row 154 -> 156).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 120..154 down to 122..156, leaving two fresh blank
# rows at 120:121 (inherits the date-format style from the row above, same
# as a manual Excel "Insert Sheet Rows").
$ws.Rows("120:121").Insert()

# New row 120 - Especial, 2021-11-11 (serial 44511)
$ws.Cells.Item(120, 1).Value = 4
$ws.Cells.Item(120, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(120, 3).Value = "Los Lagos"
$ws.Cells.Item(120, 4).Value = 44511
$ws.Cells.Item(120, 5).Value = 10
$ws.Cells.Item(120, 6).Value = "Fruta"
$ws.Cells.Item(120, 7).Value = 100101
$ws.Cells.Item(120, 8).Value = "Berries"
$ws.Cells.Item(120, 9).Value = 100101007
$ws.Cells.Item(120, 10).Value = "Kiwi"
$ws.Cells.Item(120, 11).Value = "Hayward"
$ws.Cells.Item(120, 12).Value = "Especial"
$ws.Cells.Item(120, 13).Value = 200
$ws.Cells.Item(120, 14).Value = 23000
$ws.Cells.Item(120, 15).Value = 23000
$ws.Cells.Item(120, 16).Value = 23000
$ws.Cells.Item(120, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(120, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(120, 19).Value = 1533
$ws.Cells.Item(120, 20).Value = 15

# New row 121 - Primera, 2021-11-11 (serial 44511)
$ws.Cells.Item(121, 1).Value = 4
$ws.Cells.Item(121, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(121, 3).Value = "Los Lagos"
$ws.Cells.Item(121, 4).Value = 44511
$ws.Cells.Item(121, 5).Value = 10
$ws.Cells.Item(121, 6).Value = "Fruta"
$ws.Cells.Item(121, 7).Value = 100101
$ws.Cells.Item(121, 8).Value = "Berries"
$ws.Cells.Item(121, 9).Value = 100101007
$ws.Cells.Item(121, 10).Value = "Kiwi"
$ws.Cells.Item(121, 11).Value = "Hayward"
$ws.Cells.Item(121, 12).Value = "Primera"
$ws.Cells.Item(121, 13).Value = 500
$ws.Cells.Item(121, 14).Value = 16000
$ws.Cells.Item(121, 15).Value = 17000
$ws.Cells.Item(121, 16).Value = 16500
$ws.Cells.Item(121, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(121, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(121, 19).Value = 1100
$ws.Cells.Item(121, 20).Value = 15
